$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet view stays left-to-right (explicit, matches rightToLeft="0" in the target).
$excel.ActiveWindow.DisplayRightToLeft = $false

# Append the new vote rows (3-9) under the existing header/row-2 data.
$ws.Range("A3").Value = "bmagae200@caledonian.ac.uk"
$ws.Range("B3").Value = "bvxsay376"
$ws.Range("C3").Value = "2024-09-18T20:14:55.611Z"

$ws.Range("A4").Value = "bmagae200@caledonian.ac.uk"
$ws.Range("B4").Value = "bvxsay376"
$ws.Range("C4").Value = "2024-09-18T20:17:03.532Z"

$ws.Range("A5").Value = "bmagae200@caledonian.ac.uk"
$ws.Range("B5").Value = "vffcct4569cx"
$ws.Range("C5").Value = "2024-09-18T20:38:06.798Z"

$ws.Range("A6").Value = "bmagae200@caledonian.ac.uk"
$ws.Range("B6").Value = "bffsswgv84376"
$ws.Range("C6").Value = "2024-09-18T20:43:21.626Z"

$ws.Range("A7").Value = "fibitope@gmail.com"
$ws.Range("B7").Value = "vcvse5457c"
$ws.Range("C7").Value = "2024-09-18T20:47:44.084Z"

$ws.Range("A8").Value = "b.magae@alustudent.com"
$ws.Range("B8").Value = "hbduyquy"
$ws.Range("C8").Value = "2024-09-18T21:08:57.244Z"

$ws.Range("A9").Value = "bmagae200@caledonian.ac.uk"
$ws.Range("B9").Value = "gsdeygwiugwqi"
$ws.Range("C9").Value = "2024-09-18T21:26:43.884Z"
